# references #33 Finish first draft of CCC19 mappings.
# Re-shuffles the "table name" list in column A (the CONCATENATE formulas in
# column B recompute automatically) and appends three new tables that were
# missing from the drop-all-tables worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$tables = @(
    "concept",                                     # A1
    "vocabulary",                                  # A2 (unchanged)
    "domain",                                      # A3
    "concept_class",                                # A4
    "relationship",                                 # A5
    "drug_era",                                     # A6
    "concept_synonym",                              # A7
    "redcap_records_tmp_1",                         # A8
    "concept_ancestor",                             # A9
    "drug_strength",                                # A10
    "schema_migrations",                            # A11
    "ar_internal_metadata",                         # A12
    "concept_relationship",                         # A13
    "metadata",                                     # A14
    "cdm_source",                                   # A15
    "visit_detail",                                 # A16
    "cohort_definition",                            # A17
    "death",                                        # A18
    "observation_period",                           # A19
    "person",                                       # A20
    "procedure_occurrence",                         # A21 (unchanged)
    "source_to_concept_map",                        # A22
    "specimen",                                     # A23
    "visit_occurrence",                             # A24
    "care_site",                                    # A25
    "cohort",                                       # A26
    "cohort_attribute",                             # A27
    "condition_era",                                # A28
    "condition_occurrence",                         # A29
    "cost",                                         # A30
    "device_exposure",                              # A31
    "dose_era",                                     # A32
    "drug_exposure",                                # A33
    "fact_relationship",                            # A34
    "location",                                     # A35
    "measurement",                                  # A36
    "note",                                         # A37
    "note_nlp",                                     # A38
    "observation",                                  # A39
    "payer_plan_period",                            # A40
    "provider",                                     # A41
    "redcap_projects",                              # A42
    "redcap_data_dictionaries",                     # A43
    "redcap_variable_choice_maps",                  # A44
    "redcap_variable_choices",                      # A45
    "omop_tables",                                  # A46
    "omop_columns",                                 # A47
    "redcap_variable_child_maps",                   # A48
    "redcap_events",                                # A49
    "redcap_variables",                             # A50
    "attribute_definition",                         # A51
    "redcap_event_maps",                            # A52
    "redcap_event_map_dependents",                  # A53
    "redcap_source_links",                          # A54
    "redcap_variable_maps",                         # A55
    "redcap_derived_dates",                         # A56 (new)
    "redcap_derived_date_choice_offset_mappings",   # A57 (new)
    "redcap_records_tmp_5"                          # A58 (new)
)

for ($i = 0; $i -lt $tables.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $tables[$i]
}

# Match the author's final scroll position / selection (view scrolled down
# to show the newly-added rows, selection spanning the whole drop-statement
# column through the last populated row).
$ws.Range("B1:B58").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1 | Out-Null
